# Conversion of functions to classes
# Adds a "cost" column to the menu sheet and nudges a couple of cosmetic
# view properties (column widths / selection) to mirror the authored edit.

$wb = $excel.ActiveWorkbook

$menu = $wb.Worksheets.Item("menu")
$tables = $wb.Worksheets.Item("tables")

# --- menu!E1:E9 -------------------------------------------------------
# Header
$menu.Range("E1").Value = "cost"

# Helper pattern: stamp the cell with a Text format before assigning the
# string so Excel stores it as text (t="s") instead of inferring a number,
# then clear the formatting again so no stray style sticks to the cell.
$menu.Range("E2").NumberFormat = "@"
$menu.Range("E2").Value = "0.6"
$menu.Range("E2").ClearFormats()

$menu.Range("E3").NumberFormat = "@"
$menu.Range("E3").Value = "0.7"
$menu.Range("E3").ClearFormats()

$menu.Range("E4").NumberFormat = "@"
$menu.Range("E4").Value = "0.7"
$menu.Range("E4").ClearFormats()

$menu.Range("E5").Value = 4
$menu.Range("E6").Value = 3

$menu.Range("E7").NumberFormat = "@"
$menu.Range("E7").Value = "2.5"
$menu.Range("E7").ClearFormats()

$menu.Range("E8").NumberFormat = "@"
$menu.Range("E8").Value = "2.5"
$menu.Range("E8").ClearFormats()

$menu.Range("E9").Value = 3

# --- cosmetic touch-ups -------------------------------------------------
# Column widths nudged by the autofit that happened when the new column
# was populated.
$tables.Columns.Item(1).ColumnWidth = 9
$tables.Columns.Item(2).ColumnWidth = 11.666666666666666
$menu.Columns.Item(1).ColumnWidth = 14.5
$menu.Columns.Item(4).ColumnWidth = 9.833333333333334

# Selection left where the author's cursor ended up on each sheet; the
# workbook keeps "menu" as the active/tab-selected sheet, so select it last.
$tables.Range("C55").Select() | Out-Null
$menu.Range("E10").Select() | Out-Null
